$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2959.2222
$ws.Range("I38").Value = 106.333336
$ws.Range("J38").Value = 8665
$ws.Range("K38").Value = 319.000008
$ws.Range("L38").Value = 25995
$ws.Range("M38").Value = 52.99999200000002
$ws.Range("N38").Value = -26739

$ws.Range("H80").Value = 994.56665
$ws.Range("J80").Value = 1194.25
$ws.Range("L80").Value = 3582.75
$ws.Range("N80").Value = -5578.75

$ws.Range("H83").Value = 994.56665
$ws.Range("J83").Value = 1194.25
$ws.Range("L83").Value = 10748.25
$ws.Range("N83").Value = -20732.25

$ws.Range("H101").Value = 2083.4
$ws.Range("I101").Value = 2433.5
$ws.Range("K101").Value = 7300.5
$ws.Range("M101").Value = -5678.5

$ws.Range("H112").Value = 2761.5
$ws.Range("J112").Value = 3184.7778
$ws.Range("L112").Value = 9554.3334
$ws.Range("N112").Value = -11770.3334

$ws.Range("H125").Value = 5628.357
$ws.Range("I125").Value = 2144.111
$ws.Range("J125").Value = 11900
$ws.Range("K125").Value = 19296.999
$ws.Range("L125").Value = 107100
$ws.Range("M125").Value = -16836.999
$ws.Range("N125").Value = -112020

$ws.Range("H138").Value = 11367107
$ws.Range("I138").Value = 33336522
$ws.Range("J138").Value = 3616.6553
$ws.Range("K138").Value = 100009566
$ws.Range("L138").Value = 10849.9659
$ws.Range("M138").Value = -100004426
$ws.Range("N138").Value = -21129.9659

$ws.Range("H140").Value = 75527
$ws.Range("J140").Value = 75527
$ws.Range("L140").Value = 75527
$ws.Range("N140").Value = -85887

$ws.Range("H141").Value = 16314.286
$ws.Range("I141").Value = 25450.6
$ws.Range("J141").Value = 8008.5454
$ws.Range("K141").Value = 76351.79999999999
$ws.Range("L141").Value = 24025.6362
$ws.Range("M141").Value = -71171.79999999999
$ws.Range("N141").Value = -34385.6362


# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 698073.5600000001
$ws.Range("I2").Value = 953384.25
$ws.Range("J2").Value = 1771.7273
$ws.Range("K2").Value = 953384.25
$ws.Range("L2").Value = 1771.7273
$ws.Range("M2").Value = -953271.25
$ws.Range("N2").Value = -1997.7273

$ws.Range("H45").Value = 14241.412
$ws.Range("I45").Value = 16599.111
$ws.Range("K45").Value = 16599.111
$ws.Range("M45").Value = -16222.111

$ws.Range("H101").Value = 20989
$ws.Range("J101").Value = 20989
$ws.Range("L101").Value = 20989
$ws.Range("N101").Value = -27479

$ws.Range("H102").Value = 2299.7334
$ws.Range("I102").Value = 1710.8
$ws.Range("J102").Value = 3477.6
$ws.Range("K102").Value = 1710.8
$ws.Range("L102").Value = 3477.6
$ws.Range("M102").Value = -88.79999999999995
$ws.Range("N102").Value = -6721.6

$ws.Range("H106").Value = 45613
$ws.Range("J106").Value = 43500
$ws.Range("L106").Value = 43500
$ws.Range("N106").Value = -46024

$ws.Range("H110").Value = 4228.409
$ws.Range("I110").Value = 4451.3
$ws.Range("K110").Value = 4451.3
$ws.Range("M110").Value = -2406.3

$ws.Range("H116").Value = 698073.5600000001
$ws.Range("I116").Value = 953384.25
$ws.Range("J116").Value = 1771.7273
$ws.Range("K116").Value = 953384.25
$ws.Range("L116").Value = 1771.7273
$ws.Range("M116").Value = -951090.25
$ws.Range("N116").Value = -6359.7273

$ws.Range("H132").Value = 1846.5555
$ws.Range("I132").Value = 1788.5
$ws.Range("K132").Value = 5365.5
$ws.Range("M132").Value = -2835.5


# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 698073.5600000001
$ws.Range("I3").Value = 953384.25
$ws.Range("J3").Value = 1771.7273
$ws.Range("K3").Value = 953384.25
$ws.Range("L3").Value = 1771.7273
$ws.Range("M3").Value = -953270.25
$ws.Range("N3").Value = -1999.7273

$ws.Range("H134").Value = 2294.4883
$ws.Range("I134").Value = 2294.4883
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 6883.4649
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -4348.4649

$ws.Range("N134").ClearContents()

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3053.077
$ws.Range("I134").Value = 2777.5557
$ws.Range("J134").Value = 4824.2856
$ws.Range("K134").Value = 8332.667099999999
$ws.Range("L134").Value = 14472.8568
$ws.Range("M134").Value = -5797.667099999999
$ws.Range("N134").Value = -19542.8568

$ws.Range("H141").Value = 336782.62
$ws.Range("J141").Value = 391623.22
$ws.Range("L141").Value = 391623.22
$ws.Range("N141").Value = -401983.22


# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 3859.3333
$ws.Range("I74").Value = 3859.3333
$ws.Range("K74").Value = 11577.9999
$ws.Range("M74").Value = -10516.9999

$ws.Range("H77").Value = 3859.3333
$ws.Range("I77").Value = 3859.3333
$ws.Range("K77").Value = 34733.9997
$ws.Range("M77").Value = -29429.9997

$ws.Range("H139").Value = 3647
$ws.Range("I139").Value = 952.25
$ws.Range("K139").Value = 2856.75
$ws.Range("M139").Value = 2283.25


# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 686.8823
$ws.Range("I97").Value = 673.8
$ws.Range("K97").Value = 673.8
$ws.Range("M97").Value = -177.8

$ws.Range("H132").Value = 3047.102
$ws.Range("I132").Value = 2642.342
$ws.Range("J132").Value = 4445.364
$ws.Range("K132").Value = 7927.026
$ws.Range("L132").Value = 13336.092
$ws.Range("M132").Value = -5397.026
$ws.Range("N132").Value = -18396.092


# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2188.7368
$ws.Range("J46").Value = 2585.5
$ws.Range("L46").Value = 2585.5
$ws.Range("N46").Value = -2961.5


# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 8999
$ws.Range("I58").Value = 8999
$ws.Range("K58").Value = 8999
$ws.Range("M58").Value = -8691

$ws.Range("H104").Value = 28189.8
$ws.Range("J104").Value = 28189.8
$ws.Range("N104").Value = -35177.8

$ws.Range("H128").Value = 59598
$ws.Range("J128").Value = 59598
$ws.Range("L128").Value = 59598
$ws.Range("N128").Value = -69558

$ws.Range("H133").Value = 72799.8
$ws.Range("J133").Value = 72799.8
$ws.Range("N133").Value = -82919.8

$ws.Range("H136").Value = 1930.3396
$ws.Range("I136").Value = 1935.4359
$ws.Range("J136").Value = 1916.1428
$ws.Range("K136").Value = 5806.307699999999
$ws.Range("L136").Value = 5748.428400000001
$ws.Range("M136").Value = -3256.307699999999
$ws.Range("N136").Value = -10848.4284

$ws.Range("H139").Value = 69623.5
$ws.Range("J139").Value = 69623.5
$ws.Range("L139").Value = 69623.5
$ws.Range("N139").Value = -79903.5

